# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - row => new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    4  = 168
    5  = 1326
    6  = 18523
    7  = 387
    8  = 277
    10 = 6937
    11 = 275
    20 = 280
    25 = 41
    26 = 285
    28 = 7
    30 = 5190
    33 = 64
    35 = 78
    36 = 12185
    38 = 16
    40 = 218
    41 = 299
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (sheet4) - row => new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    4  = 168
    5  = 1326
    6  = 18523
    7  = 387
    8  = 277
    10 = 6937
    11 = 275
    20 = 280
    25 = 41
    26 = 285
    28 = 7
    30 = 5190
    35 = 64
    37 = 78
    38 = 12185
    40 = 16
    42 = 218
    43 = 299
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
